$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" '37.408.65'
Set-TextValue "D3" '2.066.70'
Set-TextValue "E3" '  -0.56%  '
Set-TextValue "E4" '  +0.02%  '
Set-TextValue "D5" '234.30'
Set-TextValue "E5" '  -0.24%  '
Set-TextValue "D6" '0.625'
Set-TextValue "E6" '  +1.45%  '
Set-TextValue "E7" '  +0.03%  '
Set-TextValue "D8" '56.99'
Set-TextValue "E8" '  -0.66%  '
Set-TextValue "D9" '0.382'
Set-TextValue "E9" '  +0.20%  '
Set-TextValue "E10" '  +0.39%  '
Set-TextValue "E11" '  +0.55%  '
Set-TextValue "D12" '2.371.68'
Set-TextValue "E12" '  -0.47%  '
Set-TextValue "D13" '14.62'
Set-TextValue "E13" '  +0.82%  '
Set-TextValue "D14" '20.69'
Set-TextValue "E14" '  -1.64%  '
Set-TextValue "D15" '0.778'
Set-TextValue "E15" '  +0.47%  '
Set-TextValue "D16" '5.14'
Set-TextValue "E16" '  -2.01%  '
Set-TextValue "D17" '2.066.00'
Set-TextValue "E17" '  -0.38%  '
Set-TextValue "D18" '37.361.85'
Set-TextValue "E18" '  -0.38%  '
Set-TextValue "D19" '6.27'
Set-TextValue "E19" '  +4.51%  '
Set-TextValue "D20" '69.46'
Set-TextValue "E20" '  +1.68%  '
Set-TextValue "E21" '  -0.02%  '
Set-TextValue "D22" '226.60'
Set-TextValue "E22" '  +1.39%  '
Set-TextValue "E23" '  +0.00%  '
Set-TextValue "E24" '  +0.84%  '
Set-TextValue "E25" '  -0.73%  '
Set-TextValue "D26" '166.88'
Set-TextValue "E26" '  +2.55%  '
Set-TextValue "E27" '  -0.89%  '
Set-TextValue "E28" '  +3.92%  '
Set-TextValue "B29" 'Kaspa'
Set-TextValue "C29" 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue "D29" '0.127'
Set-TextValue "E29" '  -3.05%  '
Set-TextValue "B30" 'EthereumClassic'
Set-TextValue "C30" 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue "D30" '19.05'
Set-TextValue "E30" '  -1.42%  '
Set-TextValue "E31" '  -0.47%  '
Set-TextValue "D32" '4.46'
Set-TextValue "E32" '  -0.22%  '
Set-TextValue "D33" '0.0617'
Set-TextValue "E33" '  -0.85%  '
Set-TextValue "E34" '  +3.97%  '
Set-TextValue "D35" '2.49'
Set-TextValue "E35" '  -2.40%  '
Set-TextValue "E36" '  +0.03%  '
Set-TextValue "E37" '  -0.06%  '
Set-TextValue "E38" '  -2.14%  '
Set-TextValue "D39" '5.68'
Set-TextValue "E39" '  -4.80%  '
Set-TextValue "E40" '  -0.22%  '
Set-TextValue "D41" '4.36'
Set-TextValue "E41" '  +0.25%  '
Set-TextValue "D42" '1.465.11'
Set-TextValue "E42" '  -0.60%  '
Set-TextValue "B43" 'Cronos'
Set-TextValue "C43" 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue "D43" '0.0940'
Set-TextValue "E43" '  -2.29%  '
Set-TextValue "B44" 'Aave'
Set-TextValue "C44" 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue "D44" '96.21'
Set-TextValue "E44" '  +1.32%  '
Set-TextValue "D45" '1.18'
Set-TextValue "E45" '  +4.15%  '
Set-TextValue "E46" '  +1.45%  '
Set-TextValue "E47" '  -1.09%  '
Set-TextValue "D48" '15.06'
Set-TextValue "E48" '  -7.05%  '
Set-TextValue "E49" '  -2.01%  '
Set-TextValue "E50" '  +0.72%  '
Set-TextValue "D51" '2.260.20'
Set-TextValue "E51" '  -0.38%  '
